# 2022FSAdates.xlsx - "Add files via upload"
#
# The active sheet (tabSelected, workbookViewId=0) is Excel's "Sheet1",
# which contains Table1 (columns: A=Date, B=Category, C=Weapon,
# D=Cancelled, E=SplitGender, F=Time, G=Rollcall, H=Name, I=Link,
# K/L=derived formula columns). The edit flips the Cancelled/SplitGender
# flags for the 2022-05-29 rows (29-35), which also ripples through the
# K/L formula results automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (U13, Weapon F): SplitGender 1 -> 0
$ws.Range("E29").Value = 0
# Row 30 (U13, Weapon E): SplitGender 1 -> 0
$ws.Range("E30").Value = 0
# Row 31 (U17, Weapon F): SplitGender 1 -> 0
$ws.Range("E31").Value = 0
# Row 32 (U17, Weapon E): Cancelled 0 -> 1
$ws.Range("D32").Value = 1
# Row 33 (Open, Weapon F): SplitGender 1 -> 0
$ws.Range("E33").Value = 0
# Row 34 (Open, Weapon E): SplitGender 1 -> 0
$ws.Range("E34").Value = 0
# Row 35 (Open, Weapon S): Cancelled 0 -> 1
$ws.Range("D35").Value = 1

# Selection moved to E34 (matches the saved sheetView's <selection>)
$ws.Range("E34").Select()
